$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-21, columns B (index 0) through K (index 9)
$data = @(
    @(-15.17667962446494, 0.4260482333292384, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, 0.8159268465957269, 1.07600428500634, -15.17667962446494, 3.776655761336899, -15.17667962446494, 1.810373386912303, -15.17667962446494, 1.626808273746307, -15.17667962446494),
    @(-15.17667962446494, 0.9602514274574608, -15.17667962446494, -15.17667962446494, -15.17667962446494, 2.453310141949538, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(3.149824474409881, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, 1.896593768678801, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(3.475608689063105, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, 2.058334218306835),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, 2.30874877511725, -15.17667962446494, 2.566004806741067, -15.17667962446494, -15.17667962446494, -15.17667962446494, 1.697449240377801),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, 2.004792962740136, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, 2.24197185944819, 1.739118077747212),
    @(-15.17667962446494, -15.17667962446494, 1.173171740536414, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, 2.024841292024685),
    @(-15.17667962446494, -15.17667962446494, 0.4061973985503547, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, 1.939471395902902, -15.17667962446494),
    @(-15.17667962446494, 0.7805358035369502, 0.6727499121674879, -15.17667962446494, -15.17667962446494, -15.17667962446494, 1.40970573210949, -15.17667962446494, 2.186745070173671, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, -15.17667962446494, 2.047256080511666, 4.321891089247958, 1.923130461435321, -15.17667962446494),
    @(-15.17667962446494, -15.17667962446494, 2.66390777253398, -15.17667962446494, -15.17667962446494, -15.17667962446494, 1.544272659126705, -15.17667962446494, -15.17667962446494, -15.17667962446494),
    @(-15.17667962446494, 2.92045234444333, 2.673464472902535, -15.17667962446494, 2.654039821110587, -15.17667962446494, 1.721980803401968, -15.17667962446494, -15.17667962446494, 2.373924387157102),
    @(-15.17667962446494, 2.501156989195326, -15.17667962446494, 2.869554473513677, -15.17667962446494, 3.104530504937423, 1.801907953684831, -15.17667962446494, -15.17667962446494, -15.17667962446494)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
